$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "gender"
$ws.Range("B1").Value = "fName"
$ws.Range("C1").Value = "lName"
$ws.Range("D1").Value = "pswd"
$ws.Range("E1").Value = "day"
$ws.Range("F1").Value = "month"
$ws.Range("G1").Value = "year"
$ws.Range("H1").Value = "comPany"
$ws.Range("I1").Value = "addr"
$ws.Range("J1").Value = "cityString"
$ws.Range("K1").Value = "stateName"
$ws.Range("L1").Value = "zip"
$ws.Range("M1").Value = "countryName"
$ws.Range("N1").Value = "mobilePhone"

# ---- Data row (row 2) ----
# Values are written in the same order the form fields were filled in
# (gender/first/last name, then company/city/state/country, then the
# month dropdown, then the address field) so the shared-string table
# comes out in the same order as the recorded session.
$ws.Range("A2").Value = "Mr"
$ws.Range("B2").Value = "naresh"
$ws.Range("C2").Value = "reddy"
$ws.Range("D2").Value = 123456
$ws.Range("E2").Value = 11
$ws.Range("G2").Value = 1990
$ws.Range("H2").Value = "ahs"
$ws.Range("J2").Value = "nyk"
$ws.Range("K2").Value = "Alaska"
$ws.Range("L2").Value = 50002
$ws.Range("M2").Value = "United States"
$ws.Range("F2").Value = "June"
$ws.Range("N2").Value = 7896541230
$ws.Range("I2").Value = "Street address, P.O. Box, Company name"

# The address placeholder text was left unfilled, so it is flagged in red
# (10pt Arial, matching the error-state styling used elsewhere in the form).
$ws.Range("I2").Font.Name = "Arial"
$ws.Range("I2").Font.Size = 10
$ws.Range("I2").Font.Color = 4207601

# ---- Extend the same row style (style used by A2:B2) down through row 6 ----
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B6").PasteSpecial(-4122)

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 17.666666666666668
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(8).ColumnWidth = 9.166666666666666
$ws.Columns.Item(9).ColumnWidth = 14.666666666666666
$ws.Columns.Item(13).ColumnWidth = 12.333333333333334
$ws.Columns.Item(14).ColumnWidth = 12.666666666666666

# ---- Page orientation ----
$ws.PageSetup.Orientation = 1

# ---- Selection ----
$null = $ws.Range("F12").Select()
